$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 12, pushing existing rows 12-18 down to 13-19
$ws.Rows.Item(12).Insert()

# Match the date number format used by the rest of column D
$ws.Range("D12").NumberFormat = $ws.Range("D13").NumberFormat

# Populate the newly inserted row 12 with this week's data point
$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(12, 3).Value = "Maule"
$ws.Range("D12").Value = 44452
$ws.Cells.Item(12, 5).Value = 7
$ws.Cells.Item(12, 6).Value = 100112026
$ws.Cells.Item(12, 7).Value = "Haba"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 200
$ws.Cells.Item(12, 11).Value = 12000
$ws.Cells.Item(12, 12).Value = 12000
$ws.Cells.Item(12, 13).Value = 12000
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 16).Value = 480
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
